$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update existing B128 value (32.1 -> 31.9)
$ws.Cells.Item(128, 2).Value2 = 31.9

# 2. Append four new quarterly rows (129-132), copying number formats from row 128
$ws.Range("A128:C128").Copy() | Out-Null
$ws.Range("A129:C132").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$newData = @(
    @(44287, 32.1, 128),
    @(44378, 32.1, 129),
    @(44470, 32.200000000000003, 130),
    @(44562, 31.9, 131)
)

$r = 129
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $r++
}

# 3. Update sheet view: selection + scroll position
$ws.Range("C130").Select()
$excel.ActiveWindow.ScrollRow = 114
$excel.ActiveWindow.ScrollColumn = 1
